$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Backlog ID renumbering (Product backlog item "8" was removed earlier,
# leaving a gap between 7 and 9; this commit closes the gap by shifting
# every ID > 7 down by one across all sheets that reference it) plus
# three user-story typo fixes on the Sprint3 sheet.
# ------------------------------------------------------------------

# --- Product sheet: renumber IDs 9..30 -> 8..29 ---
$ws = $wb.Worksheets.Item("Product")
$ws.Range("A11").Value = 8
$ws.Range("A12").Value = 9
$ws.Range("A13").Value = 10
$ws.Range("A14").Value = 11
$ws.Range("A15").Value = 12
$ws.Range("A16").Value = 13
$ws.Range("A17").Value = 14
$ws.Range("A18").Value = 15
$ws.Range("A19").Value = 16
$ws.Range("A20").Value = 17
$ws.Range("A21").Value = 18
$ws.Range("A22").Value = 19
$ws.Range("A23").Value = 20
$ws.Range("A24").Value = 21
$ws.Range("A25").Value = 22
$ws.Range("A26").Value = 23
$ws.Range("A27").Value = 24
$ws.Range("A28").Value = 25
$ws.Range("A29").Value = 26
$ws.Range("A30").Value = 27
$ws.Range("A31").Value = 28
$ws.Range("A32").Value = 29
$ws.Range("A32").Select()

# --- Sprint2 sheet: renumber IDs 11..19 -> 10..18 ---
$ws = $wb.Worksheets.Item("Sprint2")
$ws.Range("A6").Value = 10
$ws.Range("A7").Value = 11
$ws.Range("A8").Value = 12
$ws.Range("A9").Value = 13
$ws.Range("A10").Value = 14
$ws.Range("A11").Value = 16
$ws.Range("A12").Value = 17
$ws.Range("A13").Value = 18
$ws.Range("A13").Select()

# --- Sprint3 sheet: renumber IDs + fix three typo'd user stories ---
$ws = $wb.Worksheets.Item("Sprint3")
$ws.Range("A4").Value = 8
$ws.Range("A5").Value = 9
$ws.Range("B5").Value = "As a player I want to be able to take actions so that the game becomes more interesting"
$ws.Range("A6").Value = 15
$ws.Range("A7").Value = 19
$ws.Range("A8").Value = 20
$ws.Range("A9").Value = 21
$ws.Range("B9").Value = "As a player I want the game to have sounds for moving pawn , win screen and special actions to make it more interesting "
$ws.Range("A10").Value = 22
$ws.Range("B10").Value = "As a player I want pawns to be displayed and move smoothly so that it makes the game more appealing"
$ws.Range("A11").Value = 23
$ws.Range("A12").Value = 24
$ws.Range("A13").Value = 25
$ws.Range("A14").Value = 26
$ws.Range("A15").Value = 27
$ws.Range("A16").Value = 28
$ws.Range("A17").Value = 29

# Sprint3 ends up the active/selected sheet+cell, matching the saved view.
$ws.Activate()
$ws.Range("A17").Select()
